$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily data rows appended after row 196 (12-10-2021)
$data = @(
    @{ Row = 197; Date = "13-10-2021"; TCM = 119.62; TCM5 = 195.21; TCMX = 105.39 },
    @{ Row = 198; Date = "14-10-2021"; TCM = 118.89; TCM5 = 193.93; TCMX = 104.82 },
    @{ Row = 199; Date = "15-10-2021"; TCM = 118.5;  TCM5 = 193.12; TCMX = 104.49 },
    @{ Row = 200; Date = "18-10-2021"; TCM = 120;    TCM5 = 195.29; TCMX = 105.84 }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Range("A$r").Value = $item.Date
    $ws.Range("B$r").Value = $item.TCM
    $ws.Range("C$r").Value = $item.TCM5
    $ws.Range("D$r").Value = $item.TCMX
}
